$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.178731441497803
$ws.Range("B1").Value = 2.416174173355103
$ws.Range("D1").Value = 2.333215236663818
$ws.Range("E1").Value = 1.196381688117981
